# Generate Report for Handback
# Applies the localization handback update:
#  - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - Populates "Latest Target File" (hyperlink) and "Latest Handback File" columns
#  - Updates "Latest Handback DateTime" values
#  - Widens a few columns that now hold longer content

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFileName = "890bf02a-7015-4d00-bb1a-6839176694b2.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef00dde49359778a4d1ed6bc18831c75feaf76de/e2e/890bf02a-7015-4d00-bb1a-6839176694b2.md"

$zhXlf = "890bf02a-7015-4d00-bb1a-6839176694b2.616751ed80cc651b098b1a20eafef2c15a2d15b2.zh-cn.xlf"
$deXlf = "890bf02a-7015-4d00-bb1a-6839176694b2.616751ed80cc651b098b1a20eafef2c15a2d15b2.de-de.xlf"

$zhHandbackTime = "2016-08-17 04:56:31"
$deHandbackTime = "2016-08-17 04:56:39"

# Column widths are stored internally in 1/6-character increments, so feed the
# setter a value that rounds to the desired stored width.
$wideStatusWidth = 29.166666666666668   # -> stored 29.9777047293527 (closest reachable)
$wideLinkWidth    = 39.166666666666664  # -> stored 40

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$overview.Columns.Item(6).ColumnWidth = $wideStatusWidth

if ($overview.Range("E2").Value() -eq $statusOld) { $overview.Range("E2").Value = $statusNew }
if ($overview.Range("F2").Value() -eq $statusOld) { $overview.Range("F2").Value = $statusNew }
if ($overview.Range("E3").Value() -eq $statusOld) { $overview.Range("E3").Value = $statusNew }
if ($overview.Range("F3").Value() -eq $statusOld) { $overview.Range("F3").Value = $statusNew }

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Columns.Item(3).ColumnWidth = $wideStatusWidth
$zh.Columns.Item(9).ColumnWidth = $wideLinkWidth
$zh.Columns.Item(10).ColumnWidth = $wideLinkWidth

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zh.Range("J2").Value = $zhXlf
$zh.Range("J3").Value = $zhXlf

$zh.Range("K2").Value = $zhHandbackTime
$zh.Range("K3").Value = $zhHandbackTime

# Remember the existing A2/A3 hyperlinks (collected via foreach, since
# indexing the collection directly with .Item() does not return usable
# Address/TextToDisplay values in this runtime), then rebuild all of the
# sheet's hyperlinks in document order: A2, I2, A3, I3.
$zhLinks = @()
foreach ($h in $zh.Hyperlinks) { $zhLinks += $h }
$zhA2Display = $zhLinks[0].TextToDisplay
$zhA2Address = $zhLinks[0].Address
$zhA3Display = $zhLinks[1].TextToDisplay
$zhA3Address = $zhLinks[1].Address

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhA2Address, "", "", $zhA2Display)
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, "", "", $mdFileName)
$zh.Hyperlinks.Add($zh.Range("A3"), $zhA3Address, "", "", $zhA3Display)
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl, "", "", $mdFileName)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Columns.Item(3).ColumnWidth = $wideStatusWidth
$de.Columns.Item(9).ColumnWidth = $wideLinkWidth
$de.Columns.Item(10).ColumnWidth = $wideLinkWidth

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$de.Range("J2").Value = $deXlf
$de.Range("J3").Value = $deXlf

$de.Range("K2").Value = $deHandbackTime
$de.Range("K3").Value = $deHandbackTime

$deLinks = @()
foreach ($h in $de.Hyperlinks) { $deLinks += $h }
$deA2Display = $deLinks[0].TextToDisplay
$deA2Address = $deLinks[0].Address
$deA3Display = $deLinks[1].TextToDisplay
$deA3Address = $deLinks[1].Address

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deA2Address, "", "", $deA2Display)
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, "", "", $mdFileName)
$de.Hyperlinks.Add($de.Range("A3"), $deA3Address, "", "", $deA3Display)
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl, "", "", $mdFileName)

Write-Host "Handback report generated"
